$wb = $excel.ActiveWorkbook

# Mapping of row -> new F value, applied identically to sheet "展览" and "全部类型"
$updates = @{
    2  = 1116
    3  = 823
    8  = 2045
    9  = 7636
    10 = 910
    11 = 428
    12 = 354
    13 = 138
    15 = 156
    16 = 7788
    18 = 1350
    23 = 312
    26 = 20
    27 = 108
    30 = 1129
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
